$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 2087.875
$ws.Range("I40").Value = 2087.875
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 2087.875
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -1912.875
# Row 100
$ws.Range("H100").Value = 2676.4211
$ws.Range("I100").Value = 2222.4614
$ws.Range("K100").Value = 2222.4614
$ws.Range("M100").Value = -1681.4614
# Row 101
$ws.Range("H101").Value = 1053.1666
$ws.Range("I101").Value = 925.4286
$ws.Range("K101").Value = 2776.2858
$ws.Range("M101").Value = -1154.2858
# Row 103
$ws.Range("H103").Value = 1359.6086
$ws.Range("I103").Value = 4154.8
$ws.Range("J103").Value = 583.1667
$ws.Range("K103").Value = 12464.4
$ws.Range("L103").Value = 1749.5001
$ws.Range("M103").Value = -11878.4
$ws.Range("N103").Value = -2921.5001
# Row 107
$ws.Range("H107").Value = 17243590
$ws.Range("I107").Value = 23810290
$ws.Range("K107").Value = 23810290
$ws.Range("M107").Value = -23808370
# Row 132
$ws.Range("H132").Value = 1453.5358
$ws.Range("I132").Value = 1107.0869
$ws.Range("J132").Value = 3047.2
$ws.Range("K132").Value = 3321.2607
$ws.Range("L132").Value = 9141.599999999999
$ws.Range("M132").Value = -791.2606999999998
$ws.Range("N132").Value = -14201.6
# Row 138
$ws.Range("H138").Value = 2839.5405
$ws.Range("I138").Value = 2148.5789
$ws.Range("J138").Value = 3568.889
$ws.Range("K138").Value = 6445.736699999999
$ws.Range("L138").Value = 10706.667
$ws.Range("M138").Value = -1305.736699999999
$ws.Range("N138").Value = -20986.667

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 97
$ws.Range("H97").Value = 1373157.4
$ws.Range("I97").Value = 1611511.8
$ws.Range("J97").Value = 2620
$ws.Range("K97").Value = 1611511.8
$ws.Range("L97").Value = 2620
$ws.Range("M97").Value = -1611015.8
$ws.Range("N97").Value = -3612
# Row 102
$ws.Range("H102").Value = 2403.9524
$ws.Range("I102").Value = 2490.7222
$ws.Range("J102").Value = 1883.3334
$ws.Range("K102").Value = 2490.7222
$ws.Range("L102").Value = 1883.3334
$ws.Range("M102").Value = -868.7222000000002
$ws.Range("N102").Value = -5127.3334
# Row 132
$ws.Range("H132").Value = 9391.157999999999
$ws.Range("I132").Value = 5244.8335
$ws.Range("K132").Value = 15734.5005
$ws.Range("M132").Value = -13204.5005

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 553.8293
$ws.Range("I94").Value = 422.44116
$ws.Range("J94").Value = 1192
$ws.Range("K94").Value = 422.44116
$ws.Range("L94").Value = 1192
$ws.Range("M94").Value = 28.55883999999998
$ws.Range("N94").Value = -2094
# Row 99
$ws.Range("H99").Value = 921.875
$ws.Range("I99").Value = 575
$ws.Range("K99").Value = 575
$ws.Range("M99").Value = 923
# Row 134
$ws.Range("H134").Value = 5664.1025
$ws.Range("I134").Value = 3523.6365
$ws.Range("K134").Value = 10570.9095
$ws.Range("M134").Value = -8035.9095

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 2825.5
$ws.Range("I16").Value = 2100
$ws.Range("K16").Value = 2100
$ws.Range("M16").Value = -1813
# Row 99
$ws.Range("H99").Value = 1395.6666
$ws.Range("I99").Value = 1394.8
$ws.Range("K99").Value = 1394.8
$ws.Range("M99").Value = 103.2
# Row 113
$ws.Range("H113").Value = 2825.5
$ws.Range("I113").Value = 2100
$ws.Range("K113").Value = 2100
$ws.Range("M113").Value = 70
# Row 122
$ws.Range("H122").Value = 2727.1765
$ws.Range("I122").Value = 2744.5625
$ws.Range("K122").Value = 8233.6875
$ws.Range("M122").Value = -5783.6875
# Row 126
$ws.Range("H126").Value = 1395.6666
$ws.Range("I126").Value = 1394.8
$ws.Range("K126").Value = 4184.4
$ws.Range("M126").Value = -1714.4
# Row 132
$ws.Range("H132").Value = 48124.31
$ws.Range("I132").Value = 3224.6
$ws.Range("K132").Value = 9673.799999999999
$ws.Range("M132").Value = -7143.799999999999
# Row 134
$ws.Range("H134").Value = 6768.75
$ws.Range("I134").Value = 7420.7334
$ws.Range("K134").Value = 22262.2002
$ws.Range("M134").Value = -19727.2002

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 339.08
$ws.Range("I5").Value = 218.66667
$ws.Range("J5").Value = 406.8125
$ws.Range("K5").Value = 656.00001
$ws.Range("L5").Value = 1220.4375
$ws.Range("M5").Value = -544.00001
$ws.Range("N5").Value = -1444.4375
# Row 12
$ws.Range("H12").Value = 830.0833
$ws.Range("I12").Value = 825.3333
$ws.Range("J12").Value = 832.93335
$ws.Range("K12").Value = 2475.9999
$ws.Range("L12").Value = 2498.80005
$ws.Range("M12").Value = -2302.9999
$ws.Range("N12").Value = -2844.80005
# Row 88
$ws.Range("H88").Value = 6908
$ws.Range("J88").Value = 12016
$ws.Range("L88").Value = 36048
$ws.Range("N88").Value = -36904
# Row 91
$ws.Range("H91").Value = 6908
$ws.Range("J91").Value = 12016
$ws.Range("L91").Value = 36048
$ws.Range("N91").Value = -39012
# Row 92
$ws.Range("H92").Value = 918.35
$ws.Range("J92").Value = 527.1875
$ws.Range("L92").Value = 1581.5625
$ws.Range("N92").Value = -4077.5625
# Row 135
$ws.Range("H135").Value = 339.08
$ws.Range("I135").Value = 218.66667
$ws.Range("J135").Value = 406.8125
$ws.Range("K135").Value = 1968.00003
$ws.Range("L135").Value = 3661.3125
$ws.Range("M135").Value = 566.9999699999998
$ws.Range("N135").Value = -8731.3125

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 2970.2856
$ws.Range("I80").Value = 2773.25
$ws.Range("J80").Value = 3233
$ws.Range("K80").Value = 2773.25
$ws.Range("L80").Value = 3233
$ws.Range("M80").Value = -1775.25
$ws.Range("N80").Value = -5229
# Row 83
$ws.Range("H83").Value = 2970.2856
$ws.Range("I83").Value = 2773.25
$ws.Range("J83").Value = 3233
$ws.Range("K83").Value = 13866.25
$ws.Range("L83").Value = 16165
$ws.Range("M83").Value = -8874.25
$ws.Range("N83").Value = -26149
# Row 98
$ws.Range("H98").Value = 50000
$ws.Range("J98").Value = 50000
$ws.Range("L98").Value = 50000
$ws.Range("N98").Value = -55990
# Row 122
$ws.Range("H122").Value = 7395.6875
$ws.Range("I122").Value = 7688.7334
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 23066.2002
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -20616.2002
$ws.Range("N122").Value = -13900
# Row 126
$ws.Range("H126").Value = 3049.111
$ws.Range("I126").Value = 2104.6667
$ws.Range("J126").Value = 4938
$ws.Range("K126").Value = 6314.000100000001
$ws.Range("L126").Value = 14814
$ws.Range("M126").Value = -3844.000100000001
$ws.Range("N126").Value = -19754
# Row 132
$ws.Range("H132").Value = 8044.61
$ws.Range("I132").Value = 5545.393
$ws.Range("K132").Value = 16636.179
$ws.Range("M132").Value = -14106.179

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 7362.75
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 7362.75
$ws.Range("K46").Value = 0
$ws.Range("L46").ClearContents()
$ws.Range("M46").Value = 7362.75
$ws.Range("N46").Value = -7738.75
# Row 68
$ws.Range("H68").Value = 2621.875
$ws.Range("I68").Value = 2830
$ws.Range("J68").Value = 1997.5
$ws.Range("K68").Value = 2830
$ws.Range("L68").Value = 1997.5
$ws.Range("M68").Value = -2081
$ws.Range("N68").Value = -3495.5
# Row 71
$ws.Range("H71").Value = 2621.875
$ws.Range("I71").Value = 2830
$ws.Range("J71").Value = 1997.5
$ws.Range("K71").Value = 14150
$ws.Range("L71").Value = 9987.5
$ws.Range("M71").Value = -10406
$ws.Range("N71").Value = -17475.5
# Row 82
$ws.Range("H82").Value = 1086.5769
$ws.Range("J82").Value = 1556.1111
$ws.Range("L82").Value = 1556.1111
$ws.Range("N82").Value = -2278.1111
# Row 85
$ws.Range("H85").Value = 1086.5769
$ws.Range("J85").Value = 1556.1111
$ws.Range("L85").Value = 1556.1111
$ws.Range("N85").Value = -4052.1111
# Row 93
$ws.Range("H93").Value = 2345.8462
$ws.Range("I93").Value = 4888.8887
$ws.Range("K93").Value = 4888.8887
$ws.Range("M93").Value = -3640.8887
# Row 132
$ws.Range("H132").Value = 5247.6816
$ws.Range("I132").Value = 4412
$ws.Range("J132").Value = 5725.2144
$ws.Range("K132").Value = 13236
$ws.Range("L132").Value = 17175.6432
$ws.Range("M132").Value = -10706
$ws.Range("N132").Value = -22235.6432

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 31750.334
$ws.Range("I62").Value = 69000
$ws.Range("K62").Value = 69000
$ws.Range("M62").Value = -68376
# Row 65
$ws.Range("H65").Value = 31750.334
$ws.Range("I65").Value = 69000
$ws.Range("K65").Value = 345000
$ws.Range("M65").Value = -341880
# Row 96
$ws.Range("H96").Value = 9486.333000000001
$ws.Range("J96").Value = 21339.6
$ws.Range("L96").Value = 21339.6
$ws.Range("N96").Value = -24085.6
# Row 100
$ws.Range("H100").Value = 296
$ws.Range("I100").Value = 256.5
$ws.Range("J100").Value = 375
$ws.Range("K100").Value = 513
$ws.Range("L100").Value = 750
$ws.Range("M100").Value = 28
$ws.Range("N100").Value = -1832
# Row 126
$ws.Range("H126").Value = 3837
$ws.Range("I126").Value = 3837
$ws.Range("K126").Value = 11511
$ws.Range("M126").Value = -9041
# Row 132
$ws.Range("H132").Value = 5407.449
$ws.Range("I132").Value = 3582.9312
$ws.Range("J132").Value = 8053
$ws.Range("K132").Value = 10748.7936
$ws.Range("L132").Value = 24159
$ws.Range("M132").Value = -8218.793600000001
$ws.Range("N132").Value = -29219
